$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" header in column H, matching the header style used by B1:G1
$ws.Range("H1").Value = "Label"
[void]$ws.Range("G1").Copy()
[void]$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New Label column values (0/1 flags per row)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1

# Refitted values for MDD 36 (row 4) and MDD 14 (row 7) in the first (100 iterations) block
$ws.Range("D4").Value = 0.3511702316438971
$ws.Range("E4").Value = 0.6488297683561028
$ws.Range("D7").Value = 0.5522280388743558
$ws.Range("E7").Value = 0.4477719611256442

[void]$ws.Range("A1").Select()
